$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "1"
$ws.Range("B2").Value = "Irina"
$ws.Range("C2").Value = "Python"
$ws.Range("D2").Value = "Bachelor"
$ws.Range("E2").Value = "1 year"

$ws.Range("A3").Value = "2"
$ws.Range("B3").Value = "Bobby"
$ws.Range("C3").Value = "C#"
$ws.Range("D3").Value = "PhD"
$ws.Range("E3").Value = "4 years"

$ws.Range("A4").Value = "3"
$ws.Range("B4").Value = "Charles"
$ws.Range("C4").Value = "C++"
$ws.Range("D4").Value = "PhD"
$ws.Range("E4").Value = "3 years"
